$d = $word.ActiveDocument

# --- 1) Merge the split "+Prettier-Code formatter: automatically format code" runs
#        (removes the stray proofErr gramStart/gramEnd wrapper around "formatter")
#        back into a single run with identical text.
$d.Content.Find.Execute(
    "+Prettier-Code formatter: automatically format code", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "+Prettier-Code formatter: automatically format code", 2)

# --- 2) Merge the split "+Material Icon Theme:  theme for file icons" runs
#        (removes the stray proofErr gramStart/gramEnd wrapper around ":  theme")
#        back into a single run with identical text.
$d.Content.Find.Execute(
    "+Material Icon Theme:  theme for file icons", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "+Material Icon Theme:  theme for file icons", 2)

# --- 3) The old "-Set up project with Create-React-App:" paragraph (which used to sit
#        right after "Material Icon Theme") is removed entirely from that spot - its
#        paragraph mark is merged away so "Material Icon Theme" is immediately followed
#        by the "npx create-react-app@..." paragraph.
$setupParaOld = $d.Paragraphs(15)
$setupParaOld.Range.Delete()

# --- 4) A brand new paragraph "+Quokka.js" is inserted right after "Material Icon Theme".
$materialIconPara = $d.Paragraphs(14)
$insertionRange = $materialIconPara.Range
$insertionRange.Collapse(0)
$insertionRange.InsertParagraphAfter()
$quokkaPara = $d.Paragraphs(15)
$quokkaPara.Range.InsertBefore("+Quokka.js")

# --- 5) Split the "npx create-react-app@<i>version</i> <i>filename</i>" paragraph right
#        before the "version" run, producing a new trailing paragraph that keeps the
#        "version"/" "/"filename" runs.
$npxPara = $d.Paragraphs(16)
$versionRange = $npxPara.Range.Duplicate
$versionRange.Find.Execute("version")
$splitPoint = $d.Range($versionRange.Start, $versionRange.Start)
$splitPoint.InsertParagraphBefore()

# --- 6) Re-purpose the original "npx create-react-app@" run (now alone in its own
#        paragraph, still the original paragraph with lastRenderedPageBreak) to instead
#        read "-Set up project with Create-React-App:" - this re-inserts that line right
#        before the (new) "npx create-react-app@ version filename" paragraph.
$setupParaNew = $d.Paragraphs(16)
$setupTextRange = $d.Range($setupParaNew.Range.Start, $setupParaNew.Range.End - 1)
$setupTextRange.Text = "-Set up project with Create-React-App:"

# --- 7) Put back the "npx create-react-app@" text as the start of the newly split-off
#        paragraph, ahead of the "version" run.
$finalNpxPara = $d.Paragraphs(17)
$finalNpxPara.Range.InsertBefore("npx create-react-app@")
